$p = $ppt.ActivePresentation

# Append a new slide at the end, using the same "Title and Content" layout
# (layout index 2) that the rest of the deck's slides use.
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Full code can be found at the following repository"
$newSlide.Shapes.Item(2).TextFrame.TextRange.Text = "https://github.com/Satyaki9207/exploring_game_data_with_mongodb"
